# "Nghỉ phép" (Leave Request) sheet: collapse the two-round
# approve/reject columns ("... Lần đầu" / "... Lần cuối") down to a
# single "Ngày Duyệt/Từ chối" + "Trạng thái" pair, dropping the old
# "Lần cuối" columns entirely and sliding "Ghi chú" left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nghỉ phép")

# Remove the now-redundant "second round" columns:
#   I (9) = "Ngày Duyệt/Từ chối Lần cuối"
#   J (10) = "Trạng thái Lần cuối"
# Deleting column 9 twice removes both, since everything right of the
# deleted column slides left one place each time.
$ws.Columns.Item(9).Delete()
$ws.Columns.Item(9).Delete()

# Simplify the remaining first-round headers (columns G/H) since there
# is now only a single round of approval.
$ws.Range("G1").Value = "Ngày Duyệt/Từ chối"
$ws.Range("H1").Value = "Trạng thái"

# Match the new target column widths (raw stored width = ColumnWidth + 5/7).
$ws.Columns.Item(7).ColumnWidth = 28 - 5/7
$ws.Columns.Item(8).ColumnWidth = 25 - 5/7
